# Apply attendance updates described by the diff:
#  - Column F (may) goes from 0 to 1 for every student row (2 through 20)
#  - Column E (april) goes from 2 to 3 for rows 2 and 9 only

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (may) for rows 2-20: 0 -> 1
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Update column E (april) for rows 2 and 9: 2 -> 3
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(9, 5).Value = 3
